$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 285 (shifts existing rows 285:400 down to 286:401)
$ws.Rows(285).Insert()

# Populate the newly inserted row 285 with the new record's data
$ws.Range("A285").Value = 7
$ws.Range("B285").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C285").Value = "Ñuble"
$ws.Range("D285").Value = 45009
$ws.Range("D285").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E285").Value = 16
$ws.Range("F285").Value = 100114013
$ws.Range("G285").Value = "Zanahoria"
$ws.Range("H285").Value = "Sin especificar"
$ws.Range("I285").Value = "Primera"
$ws.Range("J285").Value = 70
$ws.Range("K285").Value = 7000
$ws.Range("L285").Value = 7500
$ws.Range("M285").Value = 7214
$ws.Range("N285").Value = "`$/saco 20 kilos"
$ws.Range("O285").Value = "Provincia de Diguillín"
$ws.Range("P285").Value = 361
$ws.Range("Q285").Value = 20
$ws.Range("R285").Value = "Hortaliza"
